$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 73, shifting existing rows 73:82 down to 74:83
$ws.Rows.Item(73).Insert()

# Populate the new row 73 with the new weekly price observation
$ws.Cells.Item(73, 1).Value = 11
$ws.Cells.Item(73, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(73, 3).Value = "Bíobío"
$ws.Cells.Item(73, 4).Value = 45077
$ws.Cells.Item(73, 5).Value = 8
$ws.Cells.Item(73, 6).Value = 100112031
$ws.Cells.Item(73, 7).Value = "Poroto verde"
$ws.Cells.Item(73, 8).Value = "Magnum"
$ws.Cells.Item(73, 9).Value = "Primera"
$ws.Cells.Item(73, 10).Value = 100
$ws.Cells.Item(73, 11).Value = 17000
$ws.Cells.Item(73, 12).Value = 18000
$ws.Cells.Item(73, 13).Value = 17500
$ws.Cells.Item(73, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(73, 15).Value = "Región Metropolitana"
$ws.Cells.Item(73, 16).Value = 700
$ws.Cells.Item(73, 17).Value = 25
$ws.Cells.Item(73, 18).Value = "Hortaliza"
